$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.049.59"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "2.447.63"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.56"
$ws.Range("E5").Value = "  +3.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.77"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.03"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.31"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "2.831.25"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "2.449.67"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "45.945.00"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.52"
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.18"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.07"
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.72"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.91"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.27"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  +5.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.99"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.78"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.13"
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0291"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "1.959.58"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("E48").Value = "  +9.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.15"
$ws.Range("E49").Value = "  -6.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.66"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.90"
$ws.Range("E51").Value = "  +5.70%  "
